$wb = $excel.ActiveWorkbook

# --- "Blogs" sheet: add CSS Grid related resources ---
$wsBlogs = $wb.Worksheets.Item("Blogs")

$wsBlogs.Range("C8").Value = "https://mastery.games/post/grid-item-placement/"
$wsBlogs.Range("B8").Value = "How Items Flows into CSS Grid"

$wsBlogs.Range("C10").Value = "https://rachelandrew.co.uk/archives/2015/02/04/css-grid-layout-creating-complex-grids/"
$wsBlogs.Range("B10").Value = "CSS Grid Layout - Creating Complex Grid"

$wsBlogs.Range("C12").Value = "https://gridbyexample.com/examples/example21/"
$wsBlogs.Range("B12").Value = "Nested Grid"

$wsBlogs.Range("C14").Value = "https://developer.mozilla.org/en-US/docs/Web/CSS/Viewport_concepts#what_is_a_viewport"
$wsBlogs.Range("B14").Value = "ViewPort Concept"

$wsBlogs.Range("C16").Value = "https://developer.mozilla.org/en-US/docs/Mozilla/Mobile/Viewport_meta_tag"
$wsBlogs.Range("B16").Value = "ViewPort meta tag"

$wsBlogs.Range("C18").Value = "https://css-tricks.com/snippets/css/complete-guide-grid/"
$wsBlogs.Range("B18").Value = "A Complete Guide to Grid"

$wsBlogs.Range("B18").Select()

# --- "Tools" sheet: add StartBootstrap resources ---
$wsTools = $wb.Worksheets.Item("Tools")

$wsTools.Range("C6").Value = "https://github.com/StartBootstrap/startbootstrap-grayscale"
$wsTools.Range("C8").Value = "https://startbootstrap.com/"

$wsTools.Range("C8").Select()
